$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Price values (column D) are plain-text cells in the source data (e.g. "69.051.95"
# using dot thousands separators, or trailing-zero values like "1.00"). A leading
# apostrophe forces Excel to store them as text, matching the original inlineStr type,
# instead of auto-converting them to numbers and losing formatting/precision.

$ws.Range("D2").Value = '''68.946.69'
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").Value = '''3.300.99'
$ws.Range("E3").Value = '  +1.83%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''588.21'
$ws.Range("E5").Value = '  +2.09%  '

$ws.Range("D6").Value = '''185.54'
$ws.Range("E6").Value = '  +2.88%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").Value = '''0.602'
$ws.Range("E8").Value = '  +1.99%  '

$ws.Range("E9").Value = '  +4.87%  '

$ws.Range("D10").Value = '''6.71'
$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("D11").Value = '''0.425'
$ws.Range("E11").Value = '  +3.21%  '

$ws.Range("D12").Value = '''3.882.99'
$ws.Range("E12").Value = '  +2.18%  '

$ws.Range("D13").Value = '''0.138'
$ws.Range("E13").Value = '  +0.45%  '

$ws.Range("D14").Value = '''29.18'
$ws.Range("E14").Value = '  +4.79%  '

$ws.Range("D15").Value = '''68.967.75'
$ws.Range("E15").Value = '  +2.57%  '

$ws.Range("D16").Value = '''0.0000173'
$ws.Range("E16").Value = '  +3.82%  '

$ws.Range("D17").Value = '''3.283.25'
$ws.Range("E17").Value = '  +0.91%  '

$ws.Range("D18").Value = '''5.88'
$ws.Range("E18").Value = '  +1.37%  '

$ws.Range("D19").Value = '''13.72'
$ws.Range("E19").Value = '  +2.58%  '

$ws.Range("D20").Value = '''388.58'
$ws.Range("E20").Value = '  +4.33%  '

$ws.Range("D21").Value = '''7.78'
$ws.Range("E21").Value = '  +2.85%  '

$ws.Range("D22").Value = '''71.95'
$ws.Range("E22").Value = '  +1.59%  '

$ws.Range("E23").Value = '  -0.39%  '

$ws.Range("E24").Value = '  +3.66%  '

$ws.Range("D25").Value = '''0.519'
$ws.Range("E25").Value = '  +2.06%  '

$ws.Range("D26").Value = '''9.78'
$ws.Range("E26").Value = '  +2.00%  '

$ws.Range("D27").Value = '''0.187'
$ws.Range("E27").Value = '  +3.80%  '

$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.75%  '

$ws.Range("D29").Value = '''5.81'
$ws.Range("E29").Value = '  +2.74%  '

$ws.Range("D30").Value = '''2.00'
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").Value = '''23.12'
$ws.Range("E31").Value = '  +2.50%  '

$ws.Range("E32").Value = '  +4.10%  '

$ws.Range("E33").Value = '  +5.07%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("E35").Value = '  +4.20%  '

$ws.Range("D36").Value = '''163.27'
$ws.Range("E36").Value = '  +0.91%  '

$ws.Range("D37").Value = '''1.91'

$ws.Range("D38").Value = '''0.838'
$ws.Range("E38").Value = '  -1.87%  '

$ws.Range("D39").Value = '''26.62'
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").Value = '''4.62'
$ws.Range("E40").Value = '  +5.49%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '''2.63'
$ws.Range("E41").Value = '  +1.82%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''6.66'
$ws.Range("E42").Value = '  -2.04%  '

$ws.Range("D43").Value = '''41.66'
$ws.Range("E43").Value = '  +2.89%  '

$ws.Range("D44").Value = '''0.0695'
$ws.Range("E44").Value = '  +3.36%  '

$ws.Range("D45").Value = '''25.47'
$ws.Range("E45").Value = '  -0.52%  '

$ws.Range("D46").Value = '''2.638.10'
$ws.Range("E46").Value = '  -2.20%  '

$ws.Range("D47").Value = '''340.77'
$ws.Range("E47").Value = '  -5.63%  '

$ws.Range("D48").Value = '''0.0285'
$ws.Range("E48").Value = '  +2.88%  '

$ws.Range("D49").Value = '''32.39'
$ws.Range("E49").Value = '  +5.47%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''6.33'
$ws.Range("E50").Value = '  +3.74%  '

$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Value = '''0.999'
$ws.Range("E51").Value = '  +0.38%  '
